# This script updates the "PriorAuthorizationIndicator" StructureDefinition
# workbook to a newer publication snapshot:
#   - Version bump 5.0.0 -> 6.0.0
#   - Date bump to the new publication date
#   - Publisher filled in ("Alvearie Team")
#   - The duplicated "Contact / No display for ContactDetail" rows are
#     replaced by a single "Jurisdiction / United States of America" row
#   - The root Extension row's Short/Definition text in the elements table
#     is updated to the profile-specific text instead of the generic
#     "Extension" / "An Extension" placeholders

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Metadata" property/value table ---
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(3, 2).Value  = "6.0.0"
$ws1.Cells.Item(8, 2).Value  = "2022-01-21T20:46:54+00:00"
$ws1.Cells.Item(9, 2).Value  = "Alvearie Team"
$ws1.Cells.Item(10, 1).Value = "Jurisdiction"
$ws1.Cells.Item(10, 2).Value = "United States of America"

# Row 11 duplicated the "Contact" row before the edit; remove it so the
# "Description" row (previously row 12) moves up to row 11, and the sheet
# shrinks from 21 to 20 rows.
$ws1.Rows.Item(11).Delete()

# --- Sheet 2: Elements table ---
$ws2 = $wb.Worksheets.Item(2)

# Row 2 is the root "Extension" element; give it the profile-specific
# Short (column K) and Definition (column L) text.
$ws2.Cells.Item(2, 11).Value = "Prior Authorization Indicator"
$ws2.Cells.Item(2, 12).Value = "Indicates prior authorization for the claim"
